$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("06_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("07_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("07_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("07_12", 38, 37, 0, 0, 0, 0, 0, 82),
  @("07_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("08_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("08_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("08_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("08_6", 0, 0, 0, 0, 29, 29, 35, 44),
  @("09_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("09_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("09_12", 25, 18, 19, 10, 0, 0, 0, 36),
  @("09_6", 53, 0, 0, 35, 0, 0, 51, 56),
  @("10_12", 54, 0, 0, 38, 0, 0, 0, 0),
  @("10_6", 0, 0, 0, 36, 0, 0, 0, 0),
  @("10_12", 0, 0, 0, 37, 0, 0, 55, 61),
  @("10_6", 28, 0, 0, 0, 0, 0, 0, 0),
  @("11_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("11_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("11_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("11_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("12_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("12_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("12_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("12_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("13_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("13_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("13_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("13_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("14_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("14_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("14_12", 0, 0, 0, 0, 0, 0, 0, 10),
  @("14_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("15_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("15_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("15_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("15_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("16_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("16_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("16_12", 31, 25, 0, 0, 0, 0, 53, 52),
  @("16_6", 0, 0, 0, 0, 0, 0, 5, 0),
  @("17_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("17_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("17_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("17_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("18_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("18_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("18_12", 0, 0, 0, 0, 0, 0, 0, 7),
  @("18_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("19_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("19_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("19_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("19_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("20_12", 0, 0, 0, 0, 0, 0, 0, 0),
  @("20_6", 0, 0, 0, 0, 0, 0, 0, 0),
  @("20_12", 0, 0, 0, 0, 0, 0, 0, 0)
)

$startRow = 78
for ($i = 0; $i -lt $rows.Length; $i++) {
  $row = $rows[$i]
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $row[0]
  for ($c = 1; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r, $c+1).Value = $row[$c]
  }
}
